# acad_watchlist_species.xlsx -- "adding black bear and allowing non-RG obs"
#
# The watchlist is sorted alphabetically by scientific name within each
# status block. "Ursus americanus" (black bear) belongs in the "rare
# native" block between "Typha latifolia" (row 88) and "Uvularia
# sessilifolia" (old row 89). Insert a whole new row at 89, shifting
# every following row (and its formatting) down by one, then fill in
# the three columns for the new species.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 89:186 down one row, preserving each row's own formatting
# (e.g. the cell-style override that travels with "Lonicera tatarica"
# further down the sheet).
$ws.Rows.Item(89).Insert()

# New watchlist entry: scientific name / status / in.anp columns.
$ws.Range("A89").Value = "Ursus americanus"
$ws.Range("B89").Value = "rare native"
$ws.Range("C89").Value = "P"

# Leave the view roughly where the author left it (scrolled near the new
# row, with the cell below it selected).
[void]$ws.Range("A91").Select()
$excel.ActiveWindow.ScrollRow = 70
$excel.ActiveWindow.ScrollColumn = 1
